$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values would otherwise be auto-coerced to numbers by Excel
# (single-decimal-point numeric-looking strings). Force text format, write the
# values, then clear the temporary formatting so no stray number format lingers.
$textCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D15", "D20", "D21", "D22", "D24", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "42.676.62"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").Value = "2.304.17"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "316.90"
$ws.Range("E5").Value = "  -0.04%  "
$ws.Range("D6").Value = "104.06"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("D10").Value = "40.09"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("D11").Value = "0.0906"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "8.56"
$ws.Range("E12").Value = "  +3.73%  "
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "0.996"
$ws.Range("E14").Value = "  +4.08%  "
$ws.Range("D15").Value = "15.34"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").Value = "2.652.76"
$ws.Range("E16").Value = "  +0.50%  "
$ws.Range("D17").Value = "2.312.68"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").Value = "42.599.66"
$ws.Range("E18").Value = "  +1.36%  "
$ws.Range("E19").Value = "  +3.26%  "
$ws.Range("D20").Value = "0.0000106"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "13.69"
$ws.Range("E21").Value = "  +34.37%  "
$ws.Range("D22").Value = "74.00"
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("E23").Value = "  -1.56%  "
$ws.Range("D24").Value = "268.09"
$ws.Range("E24").Value = "  -3.79%  "
$ws.Range("E25").Value = "  -0.60%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("D28").Value = "2.35"
$ws.Range("E28").Value = "  -0.82%  "
$ws.Range("D29").Value = "22.60"
$ws.Range("E29").Value = "  -0.59%  "
$ws.Range("D30").Value = "6.69"
$ws.Range("E30").Value = "  +15.88%  "
$ws.Range("D31").Value = "38.01"
$ws.Range("E31").Value = "  +6.05%  "
$ws.Range("D32").Value = "165.41"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").Value = "0.0884"
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("D34").Value = "2.69"
$ws.Range("E34").Value = "  -5.02%  "
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").Value = "0.0354"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").Value = "3.72"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("D40").Value = "2.75"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("D41").Value = "1.63"
$ws.Range("E41").Value = "  +13.52%  "
$ws.Range("D42").Value = "98.50"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").Value = "70.37"
$ws.Range("E43").Value = "  +1.85%  "
$ws.Range("D44").Value = "0.227"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "12.39"
$ws.Range("E46").Value = "  +4.17%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "117.07"
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").Value = "81.01"
$ws.Range("E48").Value = "  +5.64%  "
$ws.Range("D49").Value = "1.637.80"
$ws.Range("E49").Value = "  +3.18%  "
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").Value = "8.86"
$ws.Range("E51").Value = "  -0.34%  "

# Remove the temporary text-number-format so the cells go back to the default
# (unstyled) state, matching a plain text cell with no explicit style.
foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
